# Append a new data row (row 3) to the "Relatório" sheet, mirroring the
# shape of the existing row 2 but for a later publication date.
#
# Source row 2: A2="'202211101277" (text, literal leading apostrophe)
#               B2="ACÓRDÃO"
#               C2="07/08/2025"    (text that looks like a date)
#               D2=FALSE (boolean)
#
# New row 3:    A3="'202211101277" (same process number)
#               B3="ACÓRDÃO"
#               C3="07/08/2026"    (same publication, one year later)
#               D3=FALSE

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A3 ---------------------------------------------------------------
# The target value begins with a literal apostrophe character (it is part
# of the data, not an Excel "treat as text" quote-prefix marker). Assigning
# a *single* leading apostrophe via COM is interpreted as that quote-prefix
# marker and gets stripped, so we double it (the COM-input escape for a
# literal apostrophe) to keep it in the stored text.
$ws.Range("A3").Value = "''202211101277"
# Drop the quote-prefix cell style that the assignment above implicitly
# applied, so the cell keeps the plain/default style (matching A2).
$ws.Range("A3").Style = "Normal"

# --- B3 ---------------------------------------------------------------
$ws.Range("B3").Value = "ACÓRDÃO"

# --- C3 ---------------------------------------------------------------
# "07/08/2026" parses as a date, which would silently turn the cell into a
# numeric/date value. Force text formatting first so the literal string is
# preserved, then restore the default style so no formatting residue is
# left behind on the cell (matching C2, which also carries no explicit
# style).
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "07/08/2026"
$ws.Range("C3").Style = "Normal"

# --- D3 ---------------------------------------------------------------
$ws.Range("D3").Value = $false
